$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Remove rows 5-16 (only rows 1-4 remain) ---
$ws.Range("A5:A16").EntireRow.Delete()

# --- Row 2: new job posting ---
$ws.Range("A2").Value2 = "2025-11-28 06:28:26"
$ws.Range("B2").Value2 = "【自動化】エクセルデータ転記作業の効率化依頼"
$ws.Range("D2").Value2 = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("F2").Value2 = "https://www.lancers.jp/work/detail/5442971"
$ws.Range("G2").Value2 = 145
$ws.Range("H2").Value2 = "◆効率化,自動化"

# --- Row 3: new job posting ---
$ws.Range("A3").Value2 = "2025-11-28 06:28:26"
$ws.Range("B3").Value2 = "初回 2026年1月創業 コンサル会社のバックオフィス業務フロー設計・マニュアル化、IT導入 一括見積依頼"
$ws.Range("D3").Value2 = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("F3").Value2 = "https://www.lancers.jp/work/detail/5442904"
$ws.Range("G3").Value2 = 55
$ws.Range("H3").Value2 = "◆コンサル"

# --- Row 4: new job posting (no skill summary) ---
$ws.Range("A4").Value2 = "2025-11-28 06:28:26"
$ws.Range("B4").Value2 = "【急募】PSE認証代行をお手伝いしてくれる方募集!"
$ws.Range("D4").Value2 = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("F4").Value2 = "https://www.lancers.jp/work/detail/5443188"
$ws.Range("G4").Value2 = 10
$ws.Range("H4").ClearContents()

# --- Rebuild hyperlinks: only F2, F3, F4 should remain, pointing at the new URLs ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5442971")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5442904")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5443188")
# Adding a hyperlink stamps a fresh "Hyperlink" style xf; re-apply the named
# style so the cells reuse the original style index instead of a duplicate.
$ws.Range("F2:F4").Style = "Hyperlink"

# --- Column H width 17 -> 12 (ColumnWidth units run 0.83 under the stored char-width) ---
$ws.Columns.Item(8).ColumnWidth = 11.17
